$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Row 11 was a duplicate "Contact" / "No display for ContactDetail" row -> remove it,
# everything below shifts up by one.
$ws.Rows.Item(11).Delete()

# Version bump
$ws.Cells.Item(3, 2).Value = "6.0.0"

# Updated publication date
$ws.Cells.Item(8, 2).Value = "2022-01-21T20:46:54+00:00"

# Publisher now has a real value
$ws.Cells.Item(9, 2).Value = "Alvearie Team"

# The old duplicate "Contact" row (now row 10) becomes "Jurisdiction"
$ws.Cells.Item(10, 1).Value = "Jurisdiction"
$ws.Cells.Item(10, 2).Value = "United States of America"

# "Case Sensitive" row now carries a literal text value of "true" (not boolean TRUE).
# Writing it straight in would get auto-coerced to a Boolean cell, so stage it with a
# leading apostrophe on a scratch cell (forces text), then paste only the *value* into
# place so the destination cell's existing style/format is left untouched.
$scratch = $ws.Cells.Item(200, 26)
$scratch.Value = "'true"
$scratch.Copy()
$ws.Cells.Item(14, 2).PasteSpecial(-4163) # xlPasteValues
$scratch.Delete()
$excel.CutCopyMode = 0
